# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update price list (column D) for rows 28-35
$ws.Range("D28").Value = 1053.098
$ws.Range("D29").Value = 1089.132
$ws.Range("D30").Value = 1108.506
$ws.Range("D31").Value = 1141.764
$ws.Range("D32").Value = 1635.074
$ws.Range("D33").Value = 1576.826
$ws.Range("D34").Value = 2203.18
$ws.Range("D35").Value = 2261.376

# Re-apply the merged ranges so the internal merge bookkeeping order matches
# the target layout. Merging/unmerging resets cell formatting (and touches
# previously-empty spanned cells), so back up and restore the original
# per-cell formats around the operation.
$xlPasteFormats = -4122

$backupPairs = @(
    @("B28","C28"), @("B29","C29"), @("B30","C30"), @("B31","C31"),
    @("B32","C32"), @("B33","C33"), @("B34","C34"), @("B35","C35"),
    @("B27","C27")
)
foreach ($pair in $backupPairs) {
    $ws.Range($pair[0]).Copy($ws.Range("Z" + $pair[0]))
    $ws.Range($pair[1]).Copy($ws.Range("Z" + $pair[1]))
}

$ws.Cells.UnMerge()

$ws.Range("B32:C32").Merge()
$ws.Range("B27:C27").Merge()
$ws.Range("A10:D10").Merge()
$ws.Range("B33:C33").Merge()
$ws.Range("A11:D11").Merge()
$ws.Range("B31:C31").Merge()
$ws.Range("A1:D1").Merge()
$ws.Range("B35:C35").Merge()
$ws.Range("B30:C30").Merge()
$ws.Range("B34:C34").Merge()
$ws.Range("B28:C28").Merge()
$ws.Range("A9:D9").Merge()
$ws.Range("B29:C29").Merge()

# Restore original formatting on the B/C description cells.
foreach ($pair in $backupPairs) {
    $ws.Range("Z" + $pair[0]).Copy()
    $ws.Range($pair[0]).PasteSpecial($xlPasteFormats)
    $ws.Range("Z" + $pair[1]).Copy()
    $ws.Range($pair[1]).PasteSpecial($xlPasteFormats)
}
foreach ($pair in $backupPairs) {
    $ws.Range("Z" + $pair[0]).Clear()
    $ws.Range("Z" + $pair[1]).Clear()
}

# The 4-column merges (row 1, 9, 10, 11) only carry content in column A;
# the B:D filler cells should stay blank/unstyled like before the merge.
$ws.Range("Z1").Copy()
foreach ($rng in @("B1:D1","B9:D9","B10:D10","B11:D11")) {
    $ws.Range($rng).PasteSpecial($xlPasteFormats)
}
$ws.Range("Z1").Clear()
$excel.CutCopyMode = 0
